$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instance Info")

$ws.Range("A6").Value = "i-06f0711b81a89db89"
$ws.Range("B6").Value = "t2.micro"
$ws.Range("D6").Value = "172.31.11.26"
$ws.Range("E6").Value = "2023-11-17 08:54:32+00:00"
$ws.Range("F6").Value = "launch-wizard-2"
